$p = $ppt.ActivePresentation
$s = $p.Slides.Item(18)
$s.Shapes.Item("Rectangle 5").TextFrame.TextRange.Text = "Opcionais"
